# B6-PowerPoint.pptx edit — Thu, Jun 11, 2020  7:05:15 PM
#
# 1) Three tables (on the slides that hold the "Table_0" style) get
#    re-styled from the built-in "Table_0" style
#    {7322268C-846F-49FC-BBC0-8DF0F6E9B6F5} to the built-in table style
#    {C10AC6DE-4D06-422C-BCFC-122EC0B264B7} (Table Styles gallery pick).
#
# 2) The deck's theme ("Design") is swapped from "Integral" (Red Violet
#    colour scheme) to the default "Office Theme" colour scheme.

$p = $ppt.ActivePresentation

# --- 1. Re-style every table in the deck that still uses the old
#        "Table_0" style GUID -------------------------------------------
$oldStyle = "{7322268C-846F-49FC-BBC0-8DF0F6E9B6F5}"
$newStyle = "{C10AC6DE-4D06-422C-BCFC-122EC0B264B7}"

for ($idx = 1; $idx -le $p.Slides.Count; $idx++) {
    $s = $p.Slides.Item($idx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldStyle) {
                $tbl.ApplyStyle($newStyle)
            }
        }
    }
}

# --- 2. Switch the presentation's colour scheme from "Integral" /
#        "Red Violet" to the plain "Office Theme" colours --------------
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

$cs.Colors(1).RGB  = 0        # Dk1  -> 000000
$cs.Colors(2).RGB  = 16777215 # Lt1  -> FFFFFF
$cs.Colors(3).RGB  = 6968388  # Dk2  -> 44546A
$cs.Colors(4).RGB  = 15132391 # Lt2  -> E7E6E6
$cs.Colors(5).RGB  = 13998939 # Acc1 -> 5B9BD5
$cs.Colors(6).RGB  = 3243501  # Acc2 -> ED7D31
$cs.Colors(7).RGB  = 10855845 # Acc3 -> A5A5A5
$cs.Colors(8).RGB  = 49407    # Acc4 -> FFC000
$cs.Colors(9).RGB  = 12874308 # Acc5 -> 4472C4
$cs.Colors(10).RGB = 4697456  # Acc6 -> 70AD47
$cs.Colors(11).RGB = 12673797 # Hlink-> 0563C1
$cs.Colors(12).RGB = 7491477  # FolHl-> 954F72
